$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the two existing inline-picture runs as "do not spell/grammar
#    check" (<w:rPr><w:noProof/></w:rPr>) - mirrors Word tagging a pasted
#    screenshot run as noProof.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.NoProofing = 1
$d.Paragraphs.Item(5).Range.NoProofing = 1

# ---------------------------------------------------------------------------
# 2) Materialize a brand-new numbered-list definition (numId 2) by applying
#    default numbering to a throwaway paragraph, then reshape its levels to
#    the classic 1./a./i. hybrid-multilevel pattern, then drop the
#    throwaway paragraph again (the list definition itself stays behind in
#    numbering.xml, unused, ready for numId="2" references).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$tmpPara = $d.Paragraphs.Last
$tmpPara.Style = "List Paragraph"
$tmpPara.Range.ListFormat.ApplyNumberDefault()

$lt = $tmpPara.Range.ListFormat.ListTemplate
for ($i = 0; $i -lt $lt.ListLevels.Count; $i++) {
    $lvl = $lt.ListLevels.Item($i + 1)
    $posInCycle = $i % 3
    if ($posInCycle -eq 1) {
        $lvl.NumberStyle = 4   # wdListNumberStyleLowercaseLetter
    } elseif ($posInCycle -eq 2) {
        $lvl.NumberStyle = 2   # wdListNumberStyleLowercaseRoman
    }
}

$tmpPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Append the new "Methods to add data into snowflake table" section:
#    a blank paragraph, a Heading2 paragraph, and four numbered
#    ("List Paragraph" / numId 2) bullet items - using InsertXML so the
#    run/proofErr structure comes out byte-for-byte like real Word's
#    as-you-type spell-check markup.
# ---------------------------------------------------------------------------
$insertAt = $d.Content
$insertAt.Collapse(0)

$xmlFrag = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p>
<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>
<w:r><w:t>Methods to add data into snowflake table</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>
<w:r><w:t>Insert in worksheet</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>
<w:r><w:t>Direct upload CSV files</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>
<w:r><w:t>Variable inserts in notebook</w:t></w:r>
<w:r><w:t xml:space="preserve"> -&gt; SQL</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>
<w:r><w:t xml:space="preserve">Add form </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>streamlint</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> application</w:t></w:r>
<w:r><w:t xml:space="preserve"> -&gt; Python</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertAt.InsertXML($xmlFrag)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
